# Updated figures with the Materials 37 candidate
# Mean(alph) = (0.3509;0.1427;0.0180;0.0005;0.0204;0.1406;0.3411)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "redo IRFs..." to-do item is resolved: append the author's note explaining
# why the old IRFs were kept, and mark the row as done (strikethrough), matching
# the treatment already used for the other completed items on the sheet.
$ws.Range("A29").Value = "redo IRFs of model in app, possibly interpretations - I decided to keep the old IRFs b/c they are better for interpretation"
$ws.Range("A29").Font.Strikethrough = $true

# Two new notes added alongside the existing "done" items near the bottom of the list.
$ws.Range("B26").Value = "using materials 37 material right now"
$ws.Range("B27").Value = "edit discussion"

# Switch page orientation to portrait (reflected in the saved pageSetup).
$ws.PageSetup.Orientation = 1

# Leave the selection where the author's cursor ended up after the edits.
$ws.Range("B30").Select()
